$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$table = $ws.ListObjects.Item("Table1")

$newRow = $table.ListRows.Add()

$ws.Range("A7").Value = "Brock Rumpke"
$ws.Range("B7").Value = "5 iron"
$ws.Range("C7").Value = 208
$ws.Range("D7").Value = 207.1
$ws.Range("E7").Value = 220
$ws.Range("F7").Value = 12
$ws.Range("G7").Value = 13
$ws.Range("H7").Value = 100
$ws.Range("I7").Value = 136
$ws.Range("J7").Value = 13
$ws.Range("K7").Value = 4000
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 95
$ws.Range("N7").Formula = "=Table1[[#This Row],[Ball Speed (mph)]]/Table1[[#This Row],[Club Speed (mph)]]"
$ws.Range("N7").NumberFormat = "General"
$ws.Range("O7").Value = -2
$ws.Range("P7").Value = "Srixon ZX5"
$ws.Range("Q7").Value = "Dynamic Gold X100"
$ws.Range("R7").Value = "Full Fade"
$ws.Range("S7").Value = "Normal shot, hitting a fade"

$ws.Range("R8").Select()
